# Update "想去人数" (F column) figures across the sheets, as published by
# the gh-pages data refresh at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(3, 6).Value  = 405
$ws1.Cells.Item(4, 6).Value  = 1139
$ws1.Cells.Item(5, 6).Value  = 39
$ws1.Cells.Item(8, 6).Value  = 1064
$ws1.Cells.Item(10, 6).Value = 338
$ws1.Cells.Item(14, 6).Value = 353
$ws1.Cells.Item(15, 6).Value = 28
$ws1.Cells.Item(17, 6).Value = 467
$ws1.Cells.Item(18, 6).Value = 443
$ws1.Cells.Item(19, 6).Value = 5590
$ws1.Cells.Item(20, 6).Value = 86
$ws1.Cells.Item(21, 6).Value = 1559
$ws1.Cells.Item(22, 6).Value = 367
$ws1.Cells.Item(23, 6).Value = 4756
$ws1.Cells.Item(26, 6).Value = 1497
$ws1.Cells.Item(29, 6).Value = 648
$ws1.Cells.Item(30, 6).Value = 61

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(8, 6).Value = 97

# Sheet 3: 本地生活
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Cells.Item(4, 6).Value = 2126

# Sheet 4: 全部类型 (aggregate of the other three sheets)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(4, 6).Value  = 2126
$ws4.Cells.Item(6, 6).Value  = 405
$ws4.Cells.Item(7, 6).Value  = 1139
$ws4.Cells.Item(8, 6).Value  = 39
$ws4.Cells.Item(11, 6).Value = 1064
$ws4.Cells.Item(12, 6).Value = 338
$ws4.Cells.Item(16, 6).Value = 353
$ws4.Cells.Item(17, 6).Value = 28
$ws4.Cells.Item(22, 6).Value = 443
$ws4.Cells.Item(23, 6).Value = 5590
$ws4.Cells.Item(24, 6).Value = 86
$ws4.Cells.Item(25, 6).Value = 1559
$ws4.Cells.Item(28, 6).Value = 367
$ws4.Cells.Item(31, 6).Value = 4756
$ws4.Cells.Item(34, 6).Value = 1497
$ws4.Cells.Item(37, 6).Value = 648
$ws4.Cells.Item(38, 6).Value = 61
